# Daily attendance processing - 2025-10-20 18:28:06
# For every "Recorded By" cell (column G) whose value starts with the
# literal prefix "System, ", move "System" from the front of the list
# to the end, e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$prefix = "System, "

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -ne $null -and $text.StartsWith($prefix)) {
        $rest = $text.Substring($prefix.Length)
        $cell.Value = $rest + ", System"
    }
}
